$d = $word.ActiveDocument

$d.Content.Find.Execute("788×3=2364", $true, $false, $false, $false, $false, $true, 1, $false, "835×2=1670", 2) | Out-Null
$d.Content.Find.Execute("357×6=2142", $true, $false, $false, $false, $false, $true, 1, $false, "820×3=2460", 2) | Out-Null
$d.Content.Find.Execute("343×2=686", $true, $false, $false, $false, $false, $true, 1, $false, "986×5=4930", 2) | Out-Null
$d.Content.Find.Execute("935×6=5610", $true, $false, $false, $false, $false, $true, 1, $false, "133×7=931", 2) | Out-Null
$d.Content.Find.Execute("365×2=730", $true, $false, $false, $false, $false, $true, 1, $false, "706×6=4236", 2) | Out-Null
$d.Content.Find.Execute("393×4=1572", $true, $false, $false, $false, $false, $true, 1, $false, "928×5=4640", 2) | Out-Null
$d.Content.Find.Execute("784×6=4704", $true, $false, $false, $false, $false, $true, 1, $false, "375×2=750", 2) | Out-Null
$d.Content.Find.Execute("296×2=592", $true, $false, $false, $false, $false, $true, 1, $false, "275×4=1100", 2) | Out-Null
$d.Content.Find.Execute("692×8=5536", $true, $false, $false, $false, $false, $true, 1, $false, "522×7=3654", 2) | Out-Null
$d.Content.Find.Execute("371×3=1113", $true, $false, $false, $false, $false, $true, 1, $false, "196×6=1176", 2) | Out-Null
$d.Content.Find.Execute("260×5=1300", $true, $false, $false, $false, $false, $true, 1, $false, "687×2=1374", 2) | Out-Null
$d.Content.Find.Execute("979×5=4895", $true, $false, $false, $false, $false, $true, 1, $false, "588×6=3528", 2) | Out-Null
$d.Content.Find.Execute("954×8=7632", $true, $false, $false, $false, $false, $true, 1, $false, "530×3=1590", 2) | Out-Null
$d.Content.Find.Execute("824×5=4120", $true, $false, $false, $false, $false, $true, 1, $false, "724×5=3620", 2) | Out-Null
$d.Content.Find.Execute("863×7=6041", $true, $false, $false, $false, $false, $true, 1, $false, "740×3=2220", 2) | Out-Null
$d.Content.Find.Execute("718×2=1436", $true, $false, $false, $false, $false, $true, 1, $false, "444×2=888", 2) | Out-Null
$d.Content.Find.Execute("256×5=1280", $true, $false, $false, $false, $false, $true, 1, $false, "594×3=1782", 2) | Out-Null
$d.Content.Find.Execute("831×6=4986", $true, $false, $false, $false, $false, $true, 1, $false, "233×6=1398", 2) | Out-Null
$d.Content.Find.Execute("299×9=2691", $true, $false, $false, $false, $false, $true, 1, $false, "367×8=2936", 2) | Out-Null
$d.Content.Find.Execute("371×4=1484", $true, $false, $false, $false, $false, $true, 1, $false, "462×7=3234", 2) | Out-Null
$d.Content.Find.Execute("171×8=1368", $true, $false, $false, $false, $false, $true, 1, $false, "523×3=1569", 2) | Out-Null
$d.Content.Find.Execute("122×3=366", $true, $false, $false, $false, $false, $true, 1, $false, "790×3=2370", 2) | Out-Null
$d.Content.Find.Execute("190×6=1140", $true, $false, $false, $false, $false, $true, 1, $false, "104×5=520", 2) | Out-Null
$d.Content.Find.Execute("643×8=5144", $true, $false, $false, $false, $false, $true, 1, $false, "233×4=932", 2) | Out-Null
$d.Content.Find.Execute("336×3=1008", $true, $false, $false, $false, $false, $true, 1, $false, "104×6=624", 2) | Out-Null

Write-Host "Done applying replacements"